# Apply the edit described in the diff:
# - Insert two new data rows at row 84 (pushing existing rows 84-156 down to 86-158)
# - Populate the two new rows (84 and 85) with new price data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 84, shifting existing rows down
$ws.Rows.Item(84).Insert()
$ws.Rows.Item(84).Insert()

# --- Row 84 (new) ---
$ws.Cells.Item(84,1).Value = 6
$ws.Cells.Item(84,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(84,3).Value = "Metropolitana"
$ws.Cells.Item(84,4).Value2 = 44574
$ws.Cells.Item(84,5).Value = 13
$ws.Cells.Item(84,6).Value = 100112001
$ws.Cells.Item(84,7).Value = "Berenjena"
$ws.Cells.Item(84,8).Value = "Sin especificar"
$ws.Cells.Item(84,9).Value = "Primera"
$ws.Cells.Item(84,10).Value = 180
$ws.Cells.Item(84,11).Value = 7000
$ws.Cells.Item(84,12).Value = 7000
$ws.Cells.Item(84,13).Value = 7000
$ws.Cells.Item(84,14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(84,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(84,16).Value = 140
$ws.Cells.Item(84,17).Value = 50
$ws.Cells.Item(84,18).Value = "Hortaliza"

# --- Row 85 (new) ---
$ws.Cells.Item(85,1).Value = 6
$ws.Cells.Item(85,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(85,3).Value = "Metropolitana"
$ws.Cells.Item(85,4).Value2 = 44574
$ws.Cells.Item(85,5).Value = 13
$ws.Cells.Item(85,6).Value = 100112001
$ws.Cells.Item(85,7).Value = "Berenjena"
$ws.Cells.Item(85,8).Value = "Sin especificar"
$ws.Cells.Item(85,9).Value = "Primera"
$ws.Cells.Item(85,10).Value = 120
$ws.Cells.Item(85,11).Value = 9000
$ws.Cells.Item(85,12).Value = 10000
$ws.Cells.Item(85,13).Value = 9583
$ws.Cells.Item(85,14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(85,15).Value = "Región Metropolitana"
$ws.Cells.Item(85,16).Value = 160
$ws.Cells.Item(85,17).Value = 60
$ws.Cells.Item(85,18).Value = "Hortaliza"
